$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: end time correction ---
$ws.Range("C6").Value = 0.50902777777777775

# --- Row 7: new entry (02:30 - 03:20) ---
$ws.Range("B7").Value = 0.10416666666666667
$ws.Range("C7").Value = 0.1388888888888889
$ws.Range("B7:C7").NumberFormat = "h:mm"
$ws.Range("D7").Formula = "=C7-B7"
$ws.Range("E7").Value = "Selecting new variables for principal components analysis and revisting data cleaning portion"

# --- Row 8: new entry ---
$ws.Range("B8").Value = 0.31944444444444448
$ws.Range("C8").Value = 0.33333333333333331
$ws.Range("B8:C8").NumberFormat = "h:mm"
$ws.Range("D8").Formula = "=C8-B8"
$ws.Range("E8").Value = "Coding principal components scores."

# --- Row 9: new entry ---
$ws.Range("B9").Value = 0.35694444444444445
$ws.Range("C9").Value = 0.40277777777777773
$ws.Range("B9:C9").NumberFormat = "h:mm"
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("E9").Value = "Generate clusters using PC scores and quantile groupings. Visual inspection of different combinations of variables to see if we can make intuitive sense of the different groupings and combinations."

# --- Row 10: new entry ---
$ws.Range("B10").Value = 0.45833333333333331
$ws.Range("C10").Value = 0.47222222222222227
$ws.Range("B10:C10").NumberFormat = "h:mm"
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("E10").Value = "Working on a table that will show the the principal component scores across the three key dimensiosn broken down by percentiles. I think we need to break out the demographic into a race and economic dimension, though, as major metro areas like NY are occupying a stange middle ground in the current iteration. I think some things some signs of ""weak"" economic performance, like unemployment, are working against other more positive indicators, like total population and diversity. "

# --- Row 11: old Total row is removed entirely (content + formatting) ---
$ws.Range("A11:E11").Clear()

# --- Row 15: Total row moves here, now summing through D10 ---
$ws.Range("A15").Value = "Total"
$ws.Range("D15").Formula = "=SUM(D2:D10)"

# --- Update selection to E10 ---
$ws.Range("E10").Select()
